$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.Value = '''' + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '63.586.02'
Set-TextValue $ws.Range('E2') '  +0.71%  '
Set-TextValue $ws.Range('D3') '3.103.21'
Set-TextValue $ws.Range('E3') '  -0.61%  '
Set-TextValue $ws.Range('E4') '  +0.08%  '
Set-TextValue $ws.Range('D5') '583.50'
Set-TextValue $ws.Range('E5') '  -0.31%  '
Set-TextValue $ws.Range('D6') '145.39'
Set-TextValue $ws.Range('E6') '  +0.07%  '
Set-TextValue $ws.Range('E7') '  +0.07%  '
Set-TextValue $ws.Range('D8') '3.095.08'
Set-TextValue $ws.Range('E8') '  -0.70%  '
Set-TextValue $ws.Range('D9') '0.528'
Set-TextValue $ws.Range('E9') '  -0.27%  '
Set-TextValue $ws.Range('D10') '0.160'
Set-TextValue $ws.Range('E10') '  +6.84%  '
Set-TextValue $ws.Range('D11') '5.64'
Set-TextValue $ws.Range('E11') '  -2.36%  '
Set-TextValue $ws.Range('D12') '0.457'
Set-TextValue $ws.Range('E12') '  -2.51%  '
Set-TextValue $ws.Range('D13') '0.0000246'
Set-TextValue $ws.Range('E13') '  -1.05%  '
Set-TextValue $ws.Range('D14') '37.23'
Set-TextValue $ws.Range('E14') '  +4.38%  '
Set-TextValue $ws.Range('E15') '  -1.22%  '
Set-TextValue $ws.Range('D16') '3.611.33'
Set-TextValue $ws.Range('E16') '  -0.73%  '
Set-TextValue $ws.Range('D17') '63.415.15'
Set-TextValue $ws.Range('E17') '  +0.58%  '
Set-TextValue $ws.Range('D18') '7.10'
Set-TextValue $ws.Range('E18') '  -1.33%  '
Set-TextValue $ws.Range('D19') '3.097.05'
Set-TextValue $ws.Range('E19') '  -0.73%  '
Set-TextValue $ws.Range('D20') '462.71'
Set-TextValue $ws.Range('E20') '  -0.98%  '
Set-TextValue $ws.Range('D21') '14.25'
Set-TextValue $ws.Range('E21') '  +1.18%  '
Set-TextValue $ws.Range('E22') '  -0.73%  '
Set-TextValue $ws.Range('E23') '  -1.42%  '
Set-TextValue $ws.Range('D24') '81.32'
Set-TextValue $ws.Range('E24') '  -1.12%  '
Set-TextValue $ws.Range('D25') '12.89'
Set-TextValue $ws.Range('E25') '  -3.03%  '
Set-TextValue $ws.Range('D26') '2.14'
Set-TextValue $ws.Range('E26') '  -1.38%  '
Set-TextValue $ws.Range('E27') '  +0.05%  '
Set-TextValue $ws.Range('D28') '9.04'
Set-TextValue $ws.Range('E28') '  +8.99%  '
Set-TextValue $ws.Range('E30') '  -0.43%  '
Set-TextValue $ws.Range('D31') '2.20'
Set-TextValue $ws.Range('E31') '  -1.97%  '
Set-TextValue $ws.Range('D32') '6.89'
Set-TextValue $ws.Range('E32') '  +0.53%  '
Set-TextValue $ws.Range('D33') '0.111'
Set-TextValue $ws.Range('E33') '  -0.15%  '
Set-TextValue $ws.Range('D34') '26.67'
Set-TextValue $ws.Range('E34') '  -1.56%  '
Set-TextValue $ws.Range('D35') '0.0₃0854'
Set-TextValue $ws.Range('E35') '  -2.54%  '
Set-TextValue $ws.Range('D36') '3.44'
Set-TextValue $ws.Range('E36') '  +4.29%  '
Set-TextValue $ws.Range('E37') '  -1.05%  '
Set-TextValue $ws.Range('D38') '2.31'
Set-TextValue $ws.Range('E38') '  -3.61%  '
Set-TextValue $ws.Range('D39') '6.00'
Set-TextValue $ws.Range('E39') '  -1.10%  '
Set-TextValue $ws.Range('D40') '50.26'
Set-TextValue $ws.Range('E40') '  -1.43%  '
Set-TextValue $ws.Range('D41') '435.43'
Set-TextValue $ws.Range('E41') '  -0.09%  '
Set-TextValue $ws.Range('D42') '8.69'
Set-TextValue $ws.Range('E42') '  -0.53%  '
Set-TextValue $ws.Range('D43') '0.0368'
Set-TextValue $ws.Range('E43') '  -0.66%  '
Set-TextValue $ws.Range('D44') '2.881.13'
Set-TextValue $ws.Range('E44') '  -1.81%  '
Set-TextValue $ws.Range('D45') '0.270'
Set-TextValue $ws.Range('E45') '  -3.05%  '
Set-TextValue $ws.Range('E46') '  -2.78%  '
Set-TextValue $ws.Range('D47') '36.17'
Set-TextValue $ws.Range('E47') '  +3.28%  '
Set-TextValue $ws.Range('D49') '123.50'
Set-TextValue $ws.Range('E49') '  +0.06%  '
Set-TextValue $ws.Range('D50') '0.110'
Set-TextValue $ws.Range('E50') '  -1.29%  '
Set-TextValue $ws.Range('D51') '24.09'
Set-TextValue $ws.Range('E51') '  -2.01%  '
